$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.293.99'
$ws.Range('E2').Value = '  -1.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.047.22'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.60'
$ws.Range('E5').Value = '  -1.66%  '
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.97'
$ws.Range('E8').Value = '  -3.57%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.383'
$ws.Range('E9').Value = '  -2.68%  '
$ws.Range('E10').Value = '  -2.36%  '
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.73'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.349.01'
$ws.Range('E13').Value = '  -1.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.53'
$ws.Range('E14').Value = '  -3.15%  '
$ws.Range('E15').Value = '  -2.69%  '
$ws.Range('E16').Value = '  -2.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.037.56'
$ws.Range('E17').Value = '  -2.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.244.90'
$ws.Range('E18').Value = '  -1.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.66'
$ws.Range('E20').Value = '  -2.67%  '
$ws.Range('E21').Value = '  -2.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '225.61'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('E25').Value = '  -3.64%  '
$ws.Range('E26').Value = '  +3.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.53'
$ws.Range('E27').Value = '  -1.46%  '
$ws.Range('E28').Value = '  -4.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.14'
$ws.Range('E29').Value = '  -1.79%  '
$ws.Range('E30').Value = '  -4.87%  '
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.53'
$ws.Range('E32').Value = '  -3.76%  '
$ws.Range('E33').Value = '  -1.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.56'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.47'
$ws.Range('E35').Value = '  -0.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.81'
$ws.Range('E36').Value = '  -0.58%  '
$ws.Range('E37').Value = '  -3.81%  '
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.26'
$ws.Range('E39').Value = '  -2.31%  '
$ws.Range('E40').Value = '  +4.75%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '97.76'
$ws.Range('E41').Value = '  -1.66%  '
$ws.Range('B42').Value = 'Cronos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0951'
$ws.Range('E42').Value = '  -2.64%  '
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.90'
$ws.Range('E43').Value = '  +0.45%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.479.26'
$ws.Range('E44').Value = '  +2.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.19'
$ws.Range('E45').Value = '  +3.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.58'
$ws.Range('E46').Value = '  -0.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.01'
$ws.Range('E47').Value = '  -4.03%  '
$ws.Range('E48').Value = '  -3.32%  '
$ws.Range('E49').Value = '  -2.83%  '
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.235.05'
$ws.Range('E51').Value = '  -1.48%  '
